$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (record #6)
$ws.Range("B12").Value = 43877
$ws.Range("C12").Value = 0.57638888888888895
$ws.Range("D12").Value = 0.61805555555555558
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = 60
$ws.Range("G12").Value = "Kodutoo MVC"
$ws.Range("J12").Value = 1

# Row 13 (record #7)
$ws.Range("B13").Value = 43877
$ws.Range("C13").Value = 0.65277777777777779
$ws.Range("D13").Value = 0.8125
$ws.Range("E13").Value = "-"
$ws.Range("F13").Value = 230
$ws.Range("G13").Value = "Kodutoo MVC"
$ws.Range("I13").Value = "x"
$ws.Range("J13").Value = 3

# Update selection to match the final cursor position
$ws.Range("F14").Select()

$wb.Save()
